# Add 5 more example-sentence columns + synonyms/antonyms columns to the
# word-list sheet: L:S in row 1 get new headers, formatted like the
# existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "例句3 (Sentence3)",
    "例句3中文 (CN3)",
    "例句4 (Sentence4)",
    "例句4中文 (CN4)",
    "例句5 (Sentence5)",
    "例句5中文 (CN5)",
    "近义词 (Synonyms)",
    "反义词 (Antonyms)"
)
$cols = @("L", "M", "N", "O", "P", "Q", "R", "S")

# Use the existing K1 header cell as the format template (bold header font
# + medium border + centered/wrapped/RTL-reading alignment) so the new
# header cells look exactly like the existing ones.
$ws.Range("K1").Copy()

for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.PasteSpecial(-4122)
    $cell.Value = $headers[$i]
}

$excel.CutCopyMode = $false

# Match the author's resulting selection on the header row.
$null = $ws.Range("A1:S1").Select()
